$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.7137273333333334
$ws.Cells.Item(2, 8).Value = 2.141182
$ws.Cells.Item(2, 9).Value = 0.3473329658610935
$ws.Cells.Item(2, 10).Value = 0.3473329658610935
$ws.Cells.Item(2, 13).Value = 0.1483723333333333
$ws.Cells.Item(2, 14).Value = 0.445117
$ws.Cells.Item(2, 15).Value = 0.06025689221954982
$ws.Cells.Item(2, 16).Value = 0.06025689221954982
$ws.Cells.Item(2, 17).Value = 0.1058973898104444
$ws.Cells.Item(2, 18).Value = 0.953076508294
$ws.Cells.Item(2, 19).Value = 0.02092920508818849
$ws.Cells.Item(2, 20).Value = 0.02092920508818849

$ws.Cells.Item(3, 7).Value = 0.7137273333333334
$ws.Cells.Item(3, 8).Value = 2.141182
$ws.Cells.Item(3, 9).Value = 0.3473329658610935
$ws.Cells.Item(3, 10).Value = 0.3473329658610935
$ws.Cells.Item(3, 15).Value = 0.8587848174675771
$ws.Cells.Item(3, 16).Value = 0.8587848174675772
$ws.Cells.Item(3, 17).Value = 1.509255907976444
$ws.Cells.Item(3, 18).Value = 13.583303171788
$ws.Cells.Item(3, 19).Value = 0.2982842776874914
$ws.Cells.Item(3, 20).Value = 0.2982842776874914

$ws.Cells.Item(4, 7).Value = 0.7137273333333334
$ws.Cells.Item(4, 8).Value = 2.141182
$ws.Cells.Item(4, 9).Value = 0.3473329658610935
$ws.Cells.Item(4, 10).Value = 0.3473329658610935
$ws.Cells.Item(4, 13).Value = 0.199346
$ws.Cells.Item(4, 14).Value = 0.5980380000000001
$ws.Cells.Item(4, 15).Value = 0.0809582903128731
$ws.Cells.Item(4, 16).Value = 0.0809582903128731
$ws.Cells.Item(4, 17).Value = 0.1422786889906667
$ws.Cells.Item(4, 18).Value = 1.280508200916
$ws.Cells.Item(4, 19).Value = 0.02811948308541365
$ws.Cells.Item(4, 20).Value = 0.02811948308541365

$ws.Cells.Item(5, 9).Value = 0.2938237411507374
$ws.Cells.Item(5, 10).Value = 0.2938237411507374
$ws.Cells.Item(5, 13).Value = 0.1483723333333333
$ws.Cells.Item(5, 14).Value = 0.445117
$ws.Cells.Item(5, 15).Value = 0.06025689221954982
$ws.Cells.Item(5, 16).Value = 0.06025689221954982
$ws.Cells.Item(5, 17).Value = 0.08958310989877778
$ws.Cells.Item(5, 18).Value = 0.806247989089
$ws.Cells.Item(5, 19).Value = 0.01770490550206489
$ws.Cells.Item(5, 20).Value = 0.01770490550206489

$ws.Cells.Item(6, 9).Value = 0.2938237411507374
$ws.Cells.Item(6, 10).Value = 0.2938237411507374
$ws.Cells.Item(6, 15).Value = 0.8587848174675771
$ws.Cells.Item(6, 16).Value = 0.8587848174675772
$ws.Cells.Item(6, 19).Value = 0.2523313679117767
$ws.Cells.Item(6, 20).Value = 0.2523313679117767

$ws.Cells.Item(7, 9).Value = 0.2938237411507374
$ws.Cells.Item(7, 10).Value = 0.2938237411507374
$ws.Cells.Item(7, 13).Value = 0.199346
$ws.Cells.Item(7, 14).Value = 0.5980380000000001
$ws.Cells.Item(7, 15).Value = 0.0809582903128731
$ws.Cells.Item(7, 16).Value = 0.0809582903128731
$ws.Cells.Item(7, 17).Value = 0.1203595995606667
$ws.Cells.Item(7, 18).Value = 1.083236396046
$ws.Cells.Item(7, 19).Value = 0.02378746773689588
$ws.Cells.Item(7, 20).Value = 0.02378746773689588

$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.09300700000000001
$ws.Cells.Item(8, 8).Value = 0.279021
$ws.Cells.Item(8, 9).Value = 0.04526153847152096
$ws.Cells.Item(8, 10).Value = 0.04526153847152095
$ws.Cells.Item(8, 13).Value = 0.1483723333333333
$ws.Cells.Item(8, 14).Value = 0.445117
$ws.Cells.Item(8, 15).Value = 0.06025689221954982
$ws.Cells.Item(8, 16).Value = 0.06025689221954982
$ws.Cells.Item(8, 17).Value = 0.01379966560633333
$ws.Cells.Item(8, 18).Value = 0.124196990457
$ws.Cells.Item(8, 19).Value = 0.002727319645369446
$ws.Cells.Item(8, 20).Value = 0.002727319645369446

$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.09300700000000001
$ws.Cells.Item(9, 8).Value = 0.279021
$ws.Cells.Item(9, 9).Value = 0.04526153847152096
$ws.Cells.Item(9, 10).Value = 0.04526153847152095
$ws.Cells.Item(9, 15).Value = 0.8587848174675771
$ws.Cells.Item(9, 16).Value = 0.8587848174675772
$ws.Cells.Item(9, 17).Value = 0.1966736562793333
$ws.Cells.Item(9, 18).Value = 1.770062906514
$ws.Cells.Item(9, 19).Value = 0.03886992205456685
$ws.Cells.Item(9, 20).Value = 0.03886992205456685

$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.09300700000000001
$ws.Cells.Item(10, 8).Value = 0.279021
$ws.Cells.Item(10, 9).Value = 0.04526153847152096
$ws.Cells.Item(10, 10).Value = 0.04526153847152095
$ws.Cells.Item(10, 13).Value = 0.199346
$ws.Cells.Item(10, 14).Value = 0.5980380000000001
$ws.Cells.Item(10, 15).Value = 0.0809582903128731
$ws.Cells.Item(10, 16).Value = 0.0809582903128731
$ws.Cells.Item(10, 17).Value = 0.018540573422
$ws.Cells.Item(10, 18).Value = 0.166865160798
$ws.Cells.Item(10, 19).Value = 0.003664296771584669
$ws.Cells.Item(10, 20).Value = 0.003664296771584668

$ws.Cells.Item(11, 7).Value = 0.2640916666666667
$ws.Cells.Item(11, 8).Value = 0.792275
$ws.Cells.Item(11, 9).Value = 0.1285193064053396
$ws.Cells.Item(11, 10).Value = 0.1285193064053396
$ws.Cells.Item(11, 13).Value = 0.1483723333333333
$ws.Cells.Item(11, 14).Value = 0.445117
$ws.Cells.Item(11, 15).Value = 0.06025689221954982
$ws.Cells.Item(11, 16).Value = 0.06025689221954982
$ws.Cells.Item(11, 17).Value = 0.03918389679722222
$ws.Cells.Item(11, 18).Value = 0.3526550711749999
$ws.Cells.Item(11, 19).Value = 0.007744173994197846
$ws.Cells.Item(11, 20).Value = 0.007744173994197846

$ws.Cells.Item(12, 7).Value = 0.2640916666666667
$ws.Cells.Item(12, 8).Value = 0.792275
$ws.Cells.Item(12, 9).Value = 0.1285193064053396
$ws.Cells.Item(12, 10).Value = 0.1285193064053396
$ws.Cells.Item(12, 15).Value = 0.8587848174675771
$ws.Cells.Item(12, 16).Value = 0.8587848174675772
$ws.Cells.Item(12, 17).Value = 0.5584512313722221
$ws.Cells.Item(12, 18).Value = 5.026061082349999
$ws.Cells.Item(12, 19).Value = 0.1103704290923692
$ws.Cells.Item(12, 20).Value = 0.1103704290923692

$ws.Cells.Item(13, 7).Value = 0.2640916666666667
$ws.Cells.Item(13, 8).Value = 0.792275
$ws.Cells.Item(13, 9).Value = 0.1285193064053396
$ws.Cells.Item(13, 10).Value = 0.1285193064053396
$ws.Cells.Item(13, 13).Value = 0.199346
$ws.Cells.Item(13, 14).Value = 0.5980380000000001
$ws.Cells.Item(13, 15).Value = 0.0809582903128731
$ws.Cells.Item(13, 16).Value = 0.0809582903128731
$ws.Cells.Item(13, 17).Value = 0.05264561738333334
$ws.Cells.Item(13, 18).Value = 0.47381055645
$ws.Cells.Item(13, 19).Value = 0.01040470331877258
$ws.Cells.Item(13, 20).Value = 0.01040470331877258

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2110656666666667
$ws.Cells.Item(14, 8).Value = 0.633197
$ws.Cells.Item(14, 9).Value = 0.1027143848511462
$ws.Cells.Item(14, 10).Value = 0.1027143848511461
$ws.Cells.Item(14, 13).Value = 0.1483723333333333
$ws.Cells.Item(14, 14).Value = 0.445117
$ws.Cells.Item(14, 15).Value = 0.06025689221954982
$ws.Cells.Item(14, 16).Value = 0.06025689221954982
$ws.Cells.Item(14, 17).Value = 0.03131630544988889
$ws.Cells.Item(14, 18).Value = 0.281846749049
$ws.Cells.Item(14, 19).Value = 0.006189249617372874
$ws.Cells.Item(14, 20).Value = 0.006189249617372873

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2110656666666667
$ws.Cells.Item(15, 8).Value = 0.633197
$ws.Cells.Item(15, 9).Value = 0.1027143848511462
$ws.Cells.Item(15, 10).Value = 0.1027143848511461
$ws.Cells.Item(15, 15).Value = 0.8587848174675771
$ws.Cells.Item(15, 16).Value = 0.8587848174675772
$ws.Cells.Item(15, 17).Value = 0.4463218508108888
$ws.Cells.Item(15, 18).Value = 4.016896657298
$ws.Cells.Item(15, 19).Value = 0.08820955424568602
$ws.Cells.Item(15, 20).Value = 0.08820955424568602

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2110656666666667
$ws.Cells.Item(16, 8).Value = 0.633197
$ws.Cells.Item(16, 9).Value = 0.1027143848511462
$ws.Cells.Item(16, 10).Value = 0.1027143848511461
$ws.Cells.Item(16, 13).Value = 0.199346
$ws.Cells.Item(16, 14).Value = 0.5980380000000001
$ws.Cells.Item(16, 15).Value = 0.0809582903128731
$ws.Cells.Item(16, 16).Value = 0.0809582903128731
$ws.Cells.Item(16, 17).Value = 0.04207509638733334
$ws.Cells.Item(16, 18).Value = 0.3786758674860001
$ws.Cells.Item(16, 19).Value = 0.008315580988087266
$ws.Cells.Item(16, 20).Value = 0.008315580988087264

$ws.Cells.Item(17, 7).Value = 0.1692153333333334
$ws.Cells.Item(17, 8).Value = 0.507646
$ws.Cells.Item(17, 9).Value = 0.08234806326016224
$ws.Cells.Item(17, 10).Value = 0.08234806326016222
$ws.Cells.Item(17, 13).Value = 0.1483723333333333
$ws.Cells.Item(17, 14).Value = 0.445117
$ws.Cells.Item(17, 15).Value = 0.06025689221954982
$ws.Cells.Item(17, 16).Value = 0.06025689221954982
$ws.Cells.Item(17, 17).Value = 0.02510687384244445
$ws.Cells.Item(17, 18).Value = 0.225961864582
$ws.Cells.Item(17, 19).Value = 0.004962038372356266
$ws.Cells.Item(17, 20).Value = 0.004962038372356265

$ws.Cells.Item(18, 7).Value = 0.1692153333333334
$ws.Cells.Item(18, 8).Value = 0.507646
$ws.Cells.Item(18, 9).Value = 0.08234806326016224
$ws.Cells.Item(18, 10).Value = 0.08234806326016222
$ws.Cells.Item(18, 15).Value = 0.8587848174675771
$ws.Cells.Item(18, 16).Value = 0.8587848174675772
$ws.Cells.Item(18, 17).Value = 0.3578246616404444
$ws.Cells.Item(18, 18).Value = 3.220421954764
$ws.Cells.Item(18, 19).Value = 0.07071926647568692
$ws.Cells.Item(18, 20).Value = 0.07071926647568691

$ws.Cells.Item(19, 7).Value = 0.1692153333333334
$ws.Cells.Item(19, 8).Value = 0.507646
$ws.Cells.Item(19, 9).Value = 0.08234806326016224
$ws.Cells.Item(19, 10).Value = 0.08234806326016222
$ws.Cells.Item(19, 13).Value = 0.199346
$ws.Cells.Item(19, 14).Value = 0.5980380000000001
$ws.Cells.Item(19, 15).Value = 0.0809582903128731
$ws.Cells.Item(19, 16).Value = 0.0809582903128731
$ws.Cells.Item(19, 17).Value = 0.03373239983866667
$ws.Cells.Item(19, 18).Value = 0.3035915985480001
$ws.Cells.Item(19, 19).Value = 0.006666758412119054
$ws.Cells.Item(19, 20).Value = 0.006666758412119053
